# Automatische test-sync: 2025-08-05 18:17:50
# Appends new "Logs" row (Testmail #6) and updates the "Dashboard" summary
# table to reflect the new category count.

$wb = $excel.ActiveWorkbook

# --- 1) "Logs" sheet: append the new mail-log entry as row 27 ---------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 27
$logs.Cells.Item($newRow, 1).Value = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 18:17:21"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- 2) "Logs" sheet: extend conditional-formatting ranges to row 27 --
$cfCols = @("D", "G", "H", "I", "J")
foreach ($col in $cfCols) {
    $oldRange = $logs.Range($col + "2:" + $col + "26")
    $newRange = $logs.Range($col + "2:" + $col + "27")
    for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
        $fc = $oldRange.FormatConditions.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# --- 3) "Dashboard" sheet: refresh the category summary table ---------
# The new mail pushes "Inkoop / Bestellingen" to 2 hits, so the three
# single-count categories shift down a row to keep the table sorted by
# count (descending).
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(4, 2).Value = 2

$dash.Cells.Item(5, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(5, 2).Value = 1

$dash.Cells.Item(6, 1).Value = "Klacht / Probleem"
$dash.Cells.Item(6, 2).Value = 1
